# Insert a new product row ("YASMIN 0.03/3 MG 21 TABS.") into the sales
# report just above the existing "جنتيانا نقط" row (the sheet's row 20),
# pushing the rows below it -- and the totals/footer rows -- down by one,
# and roll the new quantity into the grand total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row right before the totals row (row 23). This leaves
# the merged layout of the still-in-place product rows 20-22 untouched and
# shifts the totals row (23->24) and the footer row (24->25) down, exactly
# like the source workbook growing by one data row.
$ws.Rows.Item(23).Insert()

# The freshly inserted row comes back without the source row's borders and
# shading, so copy those (and the row height) down from the row above it.
$ws.Range("A22:N22").Copy()
$ws.Range("A23:N23").PasteSpecial(-4122)
$ws.Rows.Item(23).RowHeight = $ws.Rows.Item(22).RowHeight

# Recreate the label/quantity merges for the new row -- Insert/PasteSpecial
# don't bring merged ranges along with them.
$ws.Range("B23:G23").Merge()
$ws.Range("H23:K23").Merge()
$ws.Range("L23:M23").Merge()

# Push the three pre-existing product rows down into 21-23 ...
$ws.Range("A23").Value2 = 20
$ws.Range("B23").Value2 = "كريم فاتيكا 125 مل"
$ws.Range("H23").Value2 = "2:0"
$ws.Range("L23").Value2 = 50
$ws.Range("N23").Value2 = "1:0"

$ws.Range("A22").Value2 = 19
$ws.Range("B22").Value2 = "سرنجات 5 سم"
$ws.Range("H22").Value2 = "-1:0"
$ws.Range("L22").Value2 = 2
$ws.Range("N22").Value2 = "1:0"

$ws.Range("A21").Value2 = 18
$ws.Range("B21").Value2 = "جنتيانا نقط"
$ws.Range("H21").Value2 = "4:0"
$ws.Range("L21").Value2 = 14
$ws.Range("N21").Value2 = "2:0"

# ... and drop the brand-new product into row 20.
$ws.Range("A20").Value2 = 17
$ws.Range("B20").Value2 = "YASMIN 0.03/3 MG 21 TABS."
$ws.Range("H20").Value2 = "0:0"
$ws.Range("L20").Value2 = 205
$ws.Range("N20").Value2 = "1:0"

# Roll the new quantity into the grand total (now on row 24) -- add the new
# row's amount into the existing running total the same way the source
# report recomputed it (new amount first, then the rest of the column).
$ws.Range("K24").Formula = "=SUM(L20,L4:L19,L21:L23)"
$ws.Range("K24").Value2 = $ws.Range("K24").Value2

# Match the row-height tweaks that came along with the new row in the
# source edit (new row 23 picks up the standard data-row height; the
# footer row, now 25, shrinks slightly).
$ws.Rows.Item(23).RowHeight = 25.5
$ws.Rows.Item(25).RowHeight = 16.5
